$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# NOTE: worksheet references returned by Worksheets.Item(...) track POSITION,
# not identity, in this host -- inserting/copying a sheet elsewhere in the
# workbook silently "retargets" any variable that was bound to an index at
# or after the insertion point. To stay safe we:
#   1) do every sheet-count-changing operation (Copy/rename) FIRST, always
#      looking sheets up fresh by their (still unique) name right before
#      using them, and
#   2) only AFTER the final 5-sheet layout is in place, go back and fetch
#      each sheet fresh by its final name to fill in cell content.
# ---------------------------------------------------------------------------

# Starting layout:  1=library_content, 2=requirements, 3=answers
# Target layout:    1=library_meta, 2=requirements_meta, 3=requirements_content,
#                    4=answers_meta, 5=answers_content

# Duplicate "requirements" -> "requirements_meta" (original) + "requirements_content" (copy)
$wb.Worksheets.Item("requirements").Copy($null, $wb.Worksheets.Item("requirements"))
$wb.Worksheets.Item(3).Name = "requirements_content"
$wb.Worksheets.Item("requirements").Name = "requirements_meta"

# Duplicate "answers" -> "answers_meta" (original) + "answers_content" (copy)
$wb.Worksheets.Item("answers").Copy($null, $wb.Worksheets.Item("answers"))
$wb.Worksheets.Item(5).Name = "answers_content"
$wb.Worksheets.Item("answers").Name = "answers_meta"

# Rename the original library sheet.
$wb.Worksheets.Item("library_content").Name = "library_meta"

# ---------------------------------------------------------------------------
# 1) library_meta : rewrite key/value pairs (drop framework_* + tab rows,
#    rename keys by stripping the "library_" prefix)
# ---------------------------------------------------------------------------
$wsLib = $wb.Worksheets.Item("library_meta")
$wsLib.Cells.Clear()
$wsLib.Cells.Item(1, 1).Value = "type"
$wsLib.Cells.Item(1, 2).Value = "library"
$wsLib.Cells.Item(2, 1).Value = "urn"
$wsLib.Cells.Item(2, 2).Value = "urn:intuitem:risk:library:adobe-ccf-v5"
$wsLib.Cells.Item(3, 1).Value = "version"
$wsLib.Cells.Item(3, 2).Value = "'1"
$wsLib.Cells.Item(4, 1).Value = "locale"
$wsLib.Cells.Item(4, 2).Value = "en"
$wsLib.Cells.Item(5, 1).Value = "ref_id"
$wsLib.Cells.Item(5, 2).Value = "adobe-ccf-v5"
$wsLib.Cells.Item(6, 1).Value = "name"
$wsLib.Cells.Item(6, 2).Value = "Adobe CCF v5"
$wsLib.Cells.Item(7, 1).Value = "description"
$wsLib.Cells.Item(7, 2).Value = "Adobe Common Controls Framework (CCF) version 5`nhttps://www.adobe.com/trust/compliance/adobe-ccf.html"
$wsLib.Cells.Item(8, 1).Value = "copyright"
$wsLib.Cells.Item(8, 2).Value = "Creative Commons"
$wsLib.Cells.Item(9, 1).Value = "provider"
$wsLib.Cells.Item(9, 2).Value = "Adobe"
$wsLib.Cells.Item(10, 1).Value = "packager"
$wsLib.Cells.Item(10, 2).Value = "intuitem"

# ---------------------------------------------------------------------------
# 2) requirements_meta : brand new framework-level metadata sheet
# ---------------------------------------------------------------------------
$wsReqMeta = $wb.Worksheets.Item("requirements_meta")
$wsReqMeta.Cells.Clear()
$wsReqMeta.Cells.Item(1, 1).Value = "type"
$wsReqMeta.Cells.Item(1, 2).Value = "framework"
$wsReqMeta.Cells.Item(2, 1).Value = "base_urn"
$wsReqMeta.Cells.Item(2, 2).Value = "urn:intuitem:risk:req_node:adobe-ccf-v5"
$wsReqMeta.Cells.Item(3, 1).Value = "urn"
$wsReqMeta.Cells.Item(3, 2).Value = "urn:intuitem:risk:framework:adobe-ccf-v5"
$wsReqMeta.Cells.Item(4, 1).Value = "ref_id"
$wsReqMeta.Cells.Item(4, 2).Value = "adobe-ccf-v5"
$wsReqMeta.Cells.Item(5, 1).Value = "name"
$wsReqMeta.Cells.Item(5, 2).Value = "Adobe CCF v5"
$wsReqMeta.Cells.Item(6, 1).Value = "description"
$wsReqMeta.Cells.Item(6, 2).Value = "Adobe Common Controls Framework (CCF) version 5`nhttps://www.adobe.com/trust/compliance/adobe-ccf.html"
$wsReqMeta.Cells.Item(7, 1).Value = "answers_definition"
$wsReqMeta.Cells.Item(7, 2).Value = "answers"

# ---------------------------------------------------------------------------
# 3) requirements_content : same data as the old "requirements" sheet, but the
#    25 category-header rows drop their stray empty A/C/E/F cells, keeping
#    only B (depth) and D (name).
# ---------------------------------------------------------------------------
$wsReqContent = $wb.Worksheets.Item("requirements_content")
$categoryRows = @(2, 16, 23, 29, 45, 50, 55, 71, 94, 106, 146, 155, 160, 179, 190, 201, 206, 217, 220, 238, 246, 279, 296, 306, 320)
foreach ($r in $categoryRows) {
    $wsReqContent.Cells.Item($r, 1).ClearContents()
    $wsReqContent.Cells.Item($r, 3).ClearContents()
    $wsReqContent.Cells.Item($r, 5).ClearContents()
    $wsReqContent.Cells.Item($r, 6).ClearContents()
}

# ---------------------------------------------------------------------------
# 4) answers_meta : brand new answers-level metadata sheet
# ---------------------------------------------------------------------------
$wsAnsMeta = $wb.Worksheets.Item("answers_meta")
$wsAnsMeta.Cells.Clear()
$wsAnsMeta.Cells.Item(1, 1).Value = "type"
$wsAnsMeta.Cells.Item(1, 2).Value = "answers"
$wsAnsMeta.Cells.Item(2, 1).Value = "name"
$wsAnsMeta.Cells.Item(2, 2).Value = "answers"

# ---------------------------------------------------------------------------
# answers_content keeps the original "answers" sheet content untouched.
# ---------------------------------------------------------------------------

# Restore the original active-tab selection (first sheet).
$wb.Worksheets.Item("library_meta").Activate()
